# Apply the update described by the diff:
#  - Update Taxonsorteringsordning (col B) values on existing rows 2-4
#  - Append two new data rows (5 and 6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing rows ----
$ws.Range("B2").Value = 90792
$ws.Range("B3").Value = 89058
$ws.Range("B4").Value = 90785

# Helper: force a value to be stored as literal text (so date-like
# strings such as "2023-09-26" aren't auto-converted into real dates),
# while leaving the cell's style untouched (matching the source file,
# where none of these cells carry any explicit style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---- New row 5 ----
$ws.Range("A5").Value = 112370026
$ws.Range("B5").Value = 90792
Set-TextValue $ws.Range("C5") "Ovaliderad"
Set-TextValue $ws.Range("D5") "NT"
$ws.Range("E5").Value = 4361
Set-TextValue $ws.Range("F5") "Orange taggsvamp"
Set-TextValue $ws.Range("G5") "Hydnellum aurantiacum"
Set-TextValue $ws.Range("H5") "(Batsch:Fr.) P.Karst."
Set-TextValue $ws.Range("P5") "nybodarna Österulvsås, Jmt"
$ws.Range("Q5").Value = 469718
$ws.Range("R5").Value = 7039994
$ws.Range("S5").Value = 10
Set-TextValue $ws.Range("T5") "Jämtland"
Set-TextValue $ws.Range("U5") "Krokom"
Set-TextValue $ws.Range("V5") "Jämtland"
Set-TextValue $ws.Range("W5") "Offerdal"
Set-TextValue $ws.Range("Y5") "2023-09-26"
Set-TextValue $ws.Range("AA5") "2023-09-26"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
Set-TextValue $ws.Range("AW5") "Benny Öwre"
Set-TextValue $ws.Range("AX5") "Benny Öwre"

# ---- New row 6 ----
$ws.Range("A6").Value = 112370019
$ws.Range("B6").Value = 56430
Set-TextValue $ws.Range("C6") "Ovaliderad"
Set-TextValue $ws.Range("D6") "NT"
$ws.Range("E6").Value = 100109
Set-TextValue $ws.Range("F6") "Tretåig hackspett"
Set-TextValue $ws.Range("G6") "Picoides tridactylus"
Set-TextValue $ws.Range("H6") "(Linnaeus, 1758)"
Set-TextValue $ws.Range("P6") "nybodarna Österulvsås, Jmt"
$ws.Range("Q6").Value = 469645
$ws.Range("R6").Value = 7039915
$ws.Range("S6").Value = 10
Set-TextValue $ws.Range("T6") "Jämtland"
Set-TextValue $ws.Range("U6") "Krokom"
Set-TextValue $ws.Range("V6") "Jämtland"
Set-TextValue $ws.Range("W6") "Offerdal"
Set-TextValue $ws.Range("Y6") "2023-09-26"
Set-TextValue $ws.Range("AA6") "2023-09-26"
Set-TextValue $ws.Range("AC6") "ringhack äldre"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
Set-TextValue $ws.Range("AW6") "Benny Öwre"
Set-TextValue $ws.Range("AX6") "Benny Öwre"
